# Updated sprint 106 actual and result
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 - Actual output: "...Search," -> "...Search and Export"
$ws.Range("E6").Value = "1.It shows Est.value at count, No.of items and Below par             2.Search SKU, All outlets, Status, Search and Export                              "

# Row 13 - Actual output / Result for the "In New column next to UOM..." test case
$ws.Range("F13").Value = "It displayed the Total quantity and value"
$ws.Range("G13").Value = "Pass"

# Row 14 - Actual output / Result for the Export button test case
$ws.Range("F14").Value = "Once click the Export it downloaded Excel sheet in details"
$ws.Range("G14").Value = "Pass"

# Move the selection to reflect where the author left off
$ws.Range("F15").Select()
